# Actualización automática del tracker
# Adds rows 73-78 to Sheet1 with the latest tracked results.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure column B keeps date-like strings as literal text (matches the
# existing rows, which store "fecha" as plain text, not as Excel dates).
$ws.Range("B73:B78").NumberFormat = "@"

$ws.Cells.Item(73, 1).Value = 14601393
$ws.Cells.Item(73, 2).Value = "2025-09-11"
$ws.Cells.Item(73, 3).Value = "Federico Arnaboldi"
$ws.Cells.Item(73, 4).Value = "Stefano Napolitano"
$ws.Cells.Item(73, 5).Value = "Gana Stefano Napolitano"
$ws.Cells.Item(73, 6).Value = 1.53

$ws.Cells.Item(74, 1).Value = 14601689
$ws.Cells.Item(74, 2).Value = "2025-09-11"
$ws.Cells.Item(74, 3).Value = "Radu Mihai Papoe"
$ws.Cells.Item(74, 4).Value = "Jay Clarke"
$ws.Cells.Item(74, 5).Value = "Gana Radu Mihai Papoe"
$ws.Cells.Item(74, 6).Value = 2.25

$ws.Cells.Item(75, 1).Value = 14601437
$ws.Cells.Item(75, 2).Value = "2025-09-10"
$ws.Cells.Item(75, 3).Value = "Alex Rybakov"
$ws.Cells.Item(75, 4).Value = "Trevor Svajda"
$ws.Cells.Item(75, 5).Value = "Gana Trevor Svajda"
$ws.Cells.Item(75, 6).Value = 1.57

$ws.Cells.Item(76, 1).Value = 14601435
$ws.Cells.Item(76, 2).Value = "2025-09-10"
$ws.Cells.Item(76, 3).Value = "Jack Pinnington Jones"
$ws.Cells.Item(76, 4).Value = "Quinn Vandecasteele"
$ws.Cells.Item(76, 5).Value = "Gana Quinn Vandecasteele"
$ws.Cells.Item(76, 6).Value = 5

$ws.Cells.Item(77, 1).Value = 14601432
$ws.Cells.Item(77, 2).Value = "2025-09-10"
$ws.Cells.Item(77, 3).Value = "Luca Pow"
$ws.Cells.Item(77, 4).Value = "Rafael Jodar"
$ws.Cells.Item(77, 5).Value = "Gana Luca Pow"
$ws.Cells.Item(77, 6).Value = 5.5

$ws.Cells.Item(78, 1).Value = 14633460
$ws.Cells.Item(78, 2).Value = "2025-09-10"
$ws.Cells.Item(78, 3).Value = "Leyre Romero Gormaz"
$ws.Cells.Item(78, 4).Value = "Charo Esquiva Banuls"
$ws.Cells.Item(78, 5).Value = "Gana Charo Esquiva Banuls"
$ws.Cells.Item(78, 6).Value = 5

# Columns G (resultado) and H (profit) are still pending for these matches,
# same as the most recent existing rows - copy that "blank" shape down so
# the new rows carry real (empty) cells instead of no cell at all.
$ws.Range("G72:H72").Copy($ws.Range("G73:H78"))

Write-Output "Added rows 73-78 to Sheet1"
